$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Cash Deposit for 2001-12-12"
$ws.Range("B1").Value = "Credit Deposit for 2001-12-12"
